$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New SVR parameter columns (K, L, M) with headers + first data row.
#    This is the core content change described by the commit message:
#    "Added SVR parameter loading from pred_par structure and Excel files".
# ---------------------------------------------------------------------------
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# ---------------------------------------------------------------------------
# 2. Style clean-up: several cells that previously carried a redundant
#    "no-op" cell format (identical to the default "Normal" style) are
#    reset back to the plain default style.
# ---------------------------------------------------------------------------
$plainCells = @(
  "B1","C1","H1",
  "H2",
  "A5",
  "A6",
  "A7",
  "A8",
  "A11",
  "A14",
  "D15","E15","F15","G15","I15",
  "D16","E16","F16","G16","I16",
  "D17","E17","F17","G17","I17"
)
foreach ($ref in $plainCells) {
  $ws.Range($ref).Style = "Normal"
}

# Row 7 also carried that same redundant formatting at the row level (not
# just on its one cell) - clear it there too.
$ws.Rows("7").ClearFormats()

# Rows 13 and 14 only ever held that same redundant row-level formatting
# (row 13) / a single empty, redundantly-formatted cell (row 14) - i.e. no
# real content. Round-trip them through a delete+insert cycle (net shift of
# zero) so they come back completely bare, with no leftover formatting,
# while every row below (15, 16, 17, 24, ...) keeps its original row number.
$ws.Rows("13:14").Delete()
$ws.Rows("13:14").Insert()

# ---------------------------------------------------------------------------
# 3. Selection moved to J8 (cosmetic, matches the saved sheet view state).
# ---------------------------------------------------------------------------
$ws.Range("J8").Select()
